$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text representation (avoid Excel auto-converting
# numeric-looking strings like "3.80" or "1.00" into numbers and dropping trailing zeros).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '92.084.15'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.67%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.095.95'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.74%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.45'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '610.68'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.34%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -5.31%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.092.43'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.768'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.15%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -4.02%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '91.894.22'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '33.53'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -4.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.38'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -4.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.679.86'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.45%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.089.38'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.94%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.80'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.40'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -4.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.75'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -4.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '434.60'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -5.15%  '
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'PEPE'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0000197'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -4.23%  '
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'Uniswap'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.03'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.31%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -6.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '84.98'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -4.85%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.26'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -5.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.259.13'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.84%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +4.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.125'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -17.74%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.90%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.08'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.89%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -37.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.87'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +4.75%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -12.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '25.37'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -4.40%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.90'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.49%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.87'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -5.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.84'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +7.71%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.27'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.62%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '466.52'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -5.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.431'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.27'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -4.53%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '160.16'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.34%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -5.44%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -6.09%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.34%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '43.75'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.67%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.64%  '
